$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "26.944.14"
$ws.Range("D3").Value2 = "1.671.14"
$ws.Range("E4").Value2 = "  +0.01%  "
$ws.Range("D5").Value2 = "214.79"
$ws.Range("E5").Value2 = "  +0.04%  "
$ws.Range("E6").Value2 = "  +1.48%  "
$ws.Range("E7").Value2 = "  -0.03%  "
$ws.Range("E8").Value2 = "  +0.44%  "
$ws.Range("E9").Value2 = "  +0.56%  "
$ws.Range("D10").Value2 = "20.18"
$ws.Range("E10").Value2 = "  +0.02%  "
$ws.Range("E11").Value2 = "  +1.67%  "
$ws.Range("D12").Value2 = "1.906.85"
$ws.Range("D13").Value2 = "1.665.94"
$ws.Range("E13").Value2 = "  +0.86%  "
$ws.Range("E14").Value2 = "  +0.27%  "
$ws.Range("E15").Value2 = "  +1.13%  "
$ws.Range("D16").Value2 = "65.51"
$ws.Range("E16").Value2 = "  +0.57%  "
$ws.Range("D17").Value2 = "26.942.48"
$ws.Range("E17").Value2 = "  +0.02%  "
$ws.Range("D18").Value2 = "8.04"
$ws.Range("E18").Value2 = "  +3.83%  "
$ws.Range("D19").Value2 = "233.70"
$ws.Range("E19").Value2 = "  -0.96%  "
$ws.Range("E20").Value2 = "  +0.12%  "
$ws.Range("D22").Value2 = "4.42"
$ws.Range("E22").Value2 = "  +0.28%  "
$ws.Range("D23").Value2 = "9.17"
$ws.Range("E23").Value2 = "  -1.46%  "
$ws.Range("E24").Value2 = "  -1.97%  "
$ws.Range("D25").Value2 = "145.90"
$ws.Range("E25").Value2 = "  +0.66%  "
$ws.Range("E26").Value2 = "  +0.26%  "
$ws.Range("E27").Value2 = "  +0.85%  "
$ws.Range("E28").Value2 = "  -1.38%  "
$ws.Range("E29").Value2 = "  -0.10%  "
$ws.Range("E30").Value2 = "  +0.24%  "
$ws.Range("E31").Value2 = "  +0.21%  "
$ws.Range("D32").Value2 = "3.32"
$ws.Range("E32").Value2 = "  +0.56%  "
$ws.Range("D33").Value2 = "1.459.94"
$ws.Range("E33").Value2 = "  -5.24%  "
$ws.Range("E34").Value2 = "  +1.85%  "
$ws.Range("E35").Value2 = "  +2.00%  "
$ws.Range("E36").Value2 = "  +0.19%  "
$ws.Range("D37").Value2 = "0.580"
$ws.Range("E37").Value2 = "  -0.07%  "
$ws.Range("D38").Value2 = "0.899"
$ws.Range("E38").Value2 = "  +0.48%  "
$ws.Range("E39").Value2 = "  +1.07%  "
$ws.Range("E40").Value2 = "  +13.06%  "
$ws.Range("E41").Value2 = "  -3.43%  "
$ws.Range("E42").Value2 = "  +0.01%  "
$ws.Range("E43").Value2 = "  +2.99%  "
$ws.Range("D44").Value2 = "66.32"
$ws.Range("E44").Value2 = "  +0.23%  "
$ws.Range("D45").Value2 = "1.812.07"
$ws.Range("E45").Value2 = "  +1.08%  "
$ws.Range("E46").Value2 = "  +0.86%  "
$ws.Range("E47").Value2 = "  +0.91%  "
$ws.Range("D48").Value2 = "1.53"
$ws.Range("E48").Value2 = "  +1.06%  "
$ws.Range("D49").Value2 = "0.101"
$ws.Range("E49").Value2 = "  +2.86%  "
$ws.Range("D50").Value2 = "0.0507"
$ws.Range("E50").Value2 = "  +0.51%  "
$ws.Range("D51").Value2 = "7.68"
$ws.Range("E51").Value2 = "  +0.60%  "
